$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Summary" (Worksheets.Item(1))
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")

$summary.Range("A3").Value = "Each preset category has 5 levels with per-zone values:"
$summary.Range("A5").Value = "  • 2 presets to the left (lower/slower values)"
$summary.Range("A6").Value = "  • 2 presets to the right (higher/faster values)"

$summary.Range("A14").Value = "Frequency Presets (per-zone tick intervals):"
$summary.Range("A16").Value = "  (Each zone has unique tick intervals, e.g., Throat: 2.0s → 0.1s)"

# Insert a new blank row at row 17 (pushes "Chance Presets:" block and
# everything after it down by one row) then populate it.
$summary.Rows.Item(17).Insert()
$summary.Range("A17").Value = "  Slider range: 0.1s to 5.0s in 0.1s increments"

# ---------------------------------------------------------------------------
# Sheet "Frequency" (Worksheets.Item(4))
# ---------------------------------------------------------------------------
$freq = $wb.Worksheets.Item("Frequency")

# Insert 6 new rows right after the existing data row (row 2), copy that
# row's formatting down into them, then fill in the per-zone values.
$freq.Rows.Item(3).Resize(6).Insert()
$freq.Range("A2:F2").Copy()
$freq.Range("A3:F8").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$freq.Range("A2").Value = "Throat"
$freq.Range("B2").Value = "2.0s"
$freq.Range("C2").Value = "1.0s"
$freq.Range("D2").Value = "0.5s"
$freq.Range("E2").Value = "0.3s"
$freq.Range("F2").Value = "0.1s"

$freq.Range("A3").Value = "Head"
$freq.Range("B3").Value = "2.5s"
$freq.Range("C3").Value = "1.2s"
$freq.Range("D3").Value = "0.6s"
$freq.Range("E3").Value = "0.3s"
$freq.Range("F3").Value = "0.1s"

$freq.Range("A4").Value = "Neck"
$freq.Range("B4").Value = "2.0s"
$freq.Range("C4").Value = "1.0s"
$freq.Range("D4").Value = "0.5s"
$freq.Range("E4").Value = "0.25s"
$freq.Range("F4").Value = "0.1s"

$freq.Range("A5").Value = "Torso"
$freq.Range("B5").Value = "3.0s"
$freq.Range("C5").Value = "1.5s"
$freq.Range("D5").Value = "0.8s"
$freq.Range("E5").Value = "0.4s"
$freq.Range("F5").Value = "0.2s"

$freq.Range("A6").Value = "Arm"
$freq.Range("B6").Value = "3.5s"
$freq.Range("C6").Value = "1.8s"
$freq.Range("D6").Value = "1.0s"
$freq.Range("E6").Value = "0.5s"
$freq.Range("F6").Value = "0.2s"

$freq.Range("A7").Value = "Leg"
$freq.Range("B7").Value = "3.0s"
$freq.Range("C7").Value = "1.5s"
$freq.Range("D7").Value = "0.8s"
$freq.Range("E7").Value = "0.4s"
$freq.Range("F7").Value = "0.2s"

$freq.Range("A8").Value = "Dismemberment"
$freq.Range("B8").Value = "1.5s"
$freq.Range("C8").Value = "0.8s"
$freq.Range("D8").Value = "0.4s"
$freq.Range("E8").Value = "0.2s"
$freq.Range("F8").Value = "0.1s"
